$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Book Fair", "Fiction", "Fiction at Book Fair", 103213.71, 1978, 25803.4275, 494, 403, 5),
    @("Book Fair", "Historical", "Historical at Book Fair", 114095.37, 2142, 28523.8425, 535, 411, 5),
    @("Book Fair", "Medical", "Medical at Book Fair", 102694.96, 1954, 25673.74, 488, 399, 5),
    @("Book Fair", "Philosophy", "Philosophy at Book Fair", 103770.08, 2012, 25942.52, 503, 400, 5),
    @("Book Fair", "Technical", "Technical at Book Fair", 102181.7, 1960, 25545.425, 490, 391, 5),
    @("Summer Reading Event", "Technical", "Technical at Summer Reading Event", 103275.99, 1971, 20655.198, 394, 402, 6),
    @("Summer Reading Event", "Fiction", "Fiction at Summer Reading Event", 106479.3, 2118, 21295.86, 423, 417, 6),
    @("Summer Reading Event", "Historical", "Historical at Summer Reading Event", 97844.3, 1851, 19568.86, 370, 376, 6),
    @("Summer Reading Event", "Medical", "Medical at Summer Reading Event", 102145.52, 1969, 20429.104, 393, 396, 6),
    @("Summer Reading Event", "Philosophy", "Philosophy at Summer Reading Event", 95190.61, 1874, 19038.122, 374, 396, 6),
    @("Literature Symposium", "Fiction", "Fiction at Literature Symposium", 108134.2, 2042, 54067.1, 1021, 398, 3),
    @("Literature Symposium", "Historical", "Historical at Literature Symposium", 108068.12, 2007, 54034.06, 1003, 394, 3),
    @("Literature Symposium", "Medical", "Medical at Literature Symposium", 104269.63, 2027, 52134.815, 1013, 387, 3),
    @("Literature Symposium", "Philosophy", "Philosophy at Literature Symposium", 108110.68, 2064, 54055.34, 1032, 397, 3),
    @("Literature Symposium", "Technical", "Technical at Literature Symposium", 109092.54, 2077, 54546.27, 1038, 399, 3),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
}

# Update chart series ranges to extend from row 4 to row 16
$chartObjs = $ws.ChartObjects()
for ($i = 1; $i -le $chartObjs.Count; $i++) {
    $co = $chartObjs.Item($i)
    $chart = $co.Chart
    for ($j = 1; $j -le $chart.SeriesCollection().Count; $j++) {
        $ser = $chart.SeriesCollection().Item($j)
        $f = $ser.Formula
        $f = $f.Replace('$C$2:$C$4', '$C$2:$C$16')
        $f = $f.Replace('$D$2:$D$4', '$D$2:$D$16')
        $f = $f.Replace('$E$2:$E$4', '$E$2:$E$16')
        $f = $f.Replace('$F$2:$F$4', '$F$2:$F$16')
        $f = $f.Replace('$G$2:$G$4', '$G$2:$G$16')
        $ser.Formula = $f
    }
}